# Edit slide 4 ("Technologien"):
#  - Reposition "Picture 4" (the PWA/e-commerce icon) slightly to the right/up.
#  - Remove "Picture 26" (the css/sass/scss icon) entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# PowerPoint's object model reports/accepts Shape.Left/Top in points while the
# OOXML stores EMU (1 pt = 12700 EMU). A tiny epsilon is added before the
# EMU->point division so the point->EMU round-trip on save lands on the exact
# target EMU value instead of being truncated one EMU short.
$emuPerPt = 12700
$epsilon = 0.00004

# Move "Picture 4" to its new location.
$pic4 = $s.Shapes.Item("Picture 4")
$pic4.Left = (1153632 / $emuPerPt) + $epsilon
$pic4.Top  = (2233068 / $emuPerPt) + $epsilon

# Remove "Picture 26" completely.
$pic26 = $s.Shapes.Item("Picture 26")
$pic26.Delete()
